# Applies two changes described by the diff:
#  1. Rows 14 and 15 (match rows "Indice" 13 and 14) had their match data
#     (columns F..V) swapped - only the descriptive/result/odds columns
#     moved, the leading Indice/pais/torneio/temporada/data_partida
#     columns (A..E) stayed put.
#  2. A new row 42 (Indice 41) was appended for Zurich vs Grasshoppers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the F:V contents of rows 14 and 15 -----------------------

$row14 = @{}
$row15 = @{}
for ($col = 6; $col -le 22; $col++) {
    $row14[$col] = $ws.Cells.Item(14, $col).Value2
    $row15[$col] = $ws.Cells.Item(15, $col).Value2
}
for ($col = 6; $col -le 22; $col++) {
    $ws.Cells.Item(14, $col).Value2 = $row15[$col]
    $ws.Cells.Item(15, $col).Value2 = $row14[$col]
}

# --- 2. Append new row 42 (Zurich 2 - 1 Grasshoppers) ------------------

# Copy formatting from the last existing data row (41) so the new row
# keeps the same styles (bold/centered Indice column, date format, etc.)
$ws.Range("A41:V41").Copy()
$ws.Range("A42:V42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(42, 1).Value2 = 41
$ws.Cells.Item(42, 2).Value2 = "switzerland"
$ws.Cells.Item(42, 3).Value2 = "super-league"
$ws.Cells.Item(42, 4).Value2 = "2023-2024"
$ws.Cells.Item(42, 5).Value2 = 45195.85416666666
$ws.Cells.Item(42, 6).Value2 = "Zurich"
$ws.Cells.Item(42, 7).Value2 = 2
$ws.Cells.Item(42, 8).Value2 = "Grasshoppers"
$ws.Cells.Item(42, 9).Value2 = 1
$ws.Cells.Item(42, 10).Value2 = 1.49
$ws.Cells.Item(42, 11).Value2 = "23/09/2023 19:42"
$ws.Cells.Item(42, 12).Value2 = 1.58
$ws.Cells.Item(42, 13).Value2 = "26/09/2023 20:28"
$ws.Cells.Item(42, 14).Value2 = 4.45
$ws.Cells.Item(42, 15).Value2 = "23/09/2023 19:42"
$ws.Cells.Item(42, 16).Value2 = 4.28
$ws.Cells.Item(42, 17).Value2 = "26/09/2023 20:28"
$ws.Cells.Item(42, 18).Value2 = 5.73
$ws.Cells.Item(42, 19).Value2 = "23/09/2023 19:42"
$ws.Cells.Item(42, 20).Value2 = 5.83
$ws.Cells.Item(42, 21).Value2 = "26/09/2023 20:28"
$ws.Cells.Item(42, 22).Value2 = "https://www.betexplorer.com/football/switzerland/super-league/zurich-grasshoppers/W6HyUv55/"
